$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-06-08 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-09 Monday", 2) | Out-Null
$d.Content.Find.Execute("755÷9=83, 8", $true, $false, $false, $false, $false, $true, 1, $false, "722÷7=103, 1", 2) | Out-Null
$d.Content.Find.Execute("927÷9=103, 0", $true, $false, $false, $false, $false, $true, 1, $false, "878÷8=109, 6", 2) | Out-Null
$d.Content.Find.Execute("357÷5=71, 2", $true, $false, $false, $false, $false, $true, 1, $false, "748÷9=83, 1", 2) | Out-Null
$d.Content.Find.Execute("680÷4=170, 0", $true, $false, $false, $false, $false, $true, 1, $false, "133÷2=66, 1", 2) | Out-Null
$d.Content.Find.Execute("430÷9=47, 7", $true, $false, $false, $false, $false, $true, 1, $false, "167÷9=18, 5", 2) | Out-Null
$d.Content.Find.Execute("252÷3=84, 0", $true, $false, $false, $false, $false, $true, 1, $false, "781÷2=390, 1", 2) | Out-Null
$d.Content.Find.Execute("903÷7=129, 0", $true, $false, $false, $false, $false, $true, 1, $false, "966÷2=483, 0", 2) | Out-Null
$d.Content.Find.Execute("800÷7=114, 2", $true, $false, $false, $false, $false, $true, 1, $false, "419÷6=69, 5", 2) | Out-Null
$d.Content.Find.Execute("446÷5=89, 1", $true, $false, $false, $false, $false, $true, 1, $false, "758÷2=379, 0", 2) | Out-Null
$d.Content.Find.Execute("919÷2=459, 1", $true, $false, $false, $false, $false, $true, 1, $false, "819÷9=91, 0", 2) | Out-Null
$d.Content.Find.Execute("604÷6=100, 4", $true, $false, $false, $false, $false, $true, 1, $false, "807÷5=161, 2", 2) | Out-Null
$d.Content.Find.Execute("683÷2=341, 1", $true, $false, $false, $false, $false, $true, 1, $false, "870÷3=290, 0", 2) | Out-Null
$d.Content.Find.Execute("949÷5=189, 4", $true, $false, $false, $false, $false, $true, 1, $false, "370÷3=123, 1", 2) | Out-Null
$d.Content.Find.Execute("439÷5=87, 4", $true, $false, $false, $false, $false, $true, 1, $false, "229÷6=38, 1", 2) | Out-Null
$d.Content.Find.Execute("199÷4=49, 3", $true, $false, $false, $false, $false, $true, 1, $false, "682÷7=97, 3", 2) | Out-Null
$d.Content.Find.Execute("647÷3=215, 2", $true, $false, $false, $false, $false, $true, 1, $false, "992÷7=141, 5", 2) | Out-Null
$d.Content.Find.Execute("317÷4=79, 1", $true, $false, $false, $false, $false, $true, 1, $false, "108÷6=18, 0", 2) | Out-Null
$d.Content.Find.Execute("149÷8=18, 5", $true, $false, $false, $false, $false, $true, 1, $false, "852÷8=106, 4", 2) | Out-Null
$d.Content.Find.Execute("502÷9=55, 7", $true, $false, $false, $false, $false, $true, 1, $false, "724÷8=90, 4", 2) | Out-Null
$d.Content.Find.Execute("272÷6=45, 2", $true, $false, $false, $false, $false, $true, 1, $false, "706÷9=78, 4", 2) | Out-Null
$d.Content.Find.Execute("381÷7=54, 3", $true, $false, $false, $false, $false, $true, 1, $false, "144÷2=72, 0", 2) | Out-Null
$d.Content.Find.Execute("720÷7=102, 6", $true, $false, $false, $false, $false, $true, 1, $false, "396÷2=198, 0", 2) | Out-Null
$d.Content.Find.Execute("673÷5=134, 3", $true, $false, $false, $false, $false, $true, 1, $false, "284÷7=40, 4", 2) | Out-Null
$d.Content.Find.Execute("337÷6=56, 1", $true, $false, $false, $false, $false, $true, 1, $false, "997÷3=332, 1", 2) | Out-Null
$d.Content.Find.Execute("165÷3=55, 0", $true, $false, $false, $false, $false, $true, 1, $false, "656÷3=218, 2", 2) | Out-Null
